# ---------------------------------------------------------------------------
# Rebuilds the cover-letter template body:
#  - the lone (mis-formatted) placeholder paragraph becomes a full letter
#    skeleton (date / lawyer block / RE / dates / body / signature).
#  - every paragraph gets explicit `spacing after=0` (+ line=240/auto for
#    most of them) instead of inheriting the stray sz=18 / jc=left direct
#    formatting that had been sitting on the original paragraph mark.
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# Sentinel control characters used purely as split markers while building
# the text in one shot; neither appears in real template text.
$PB = [char]1   # -> paragraph break (^p)
$TB = [char]2   # -> tab            (^t)

$first = $d.Paragraphs(1)

# Strip every bit of direct formatting (spacing/jc/rPr sz=18) the original
# lone paragraph carried, so nothing stray propagates into the paragraphs
# that get split off of it below.
$first.Style = "Normal"

$lines = @(
    "February 13, 2024",
    "",
    "lawyername",
    "lawyer_office_name",
    "lawyer_office_address1",
    "lawyer_office_address2",
    "",
    "Tel: lawyer_phone",
    "",
    "",
    "RE: ${TB}case_name",
    "",
    "Date of Crash:${TB}crash_date",
    "Date of Birth:${TB}plaintiff1_dob",
    "${TB}plaintiff2_dob",
    "",
    "doc_body",
    "signature_block"
)

$full = [string]::Join([string]$PB, $lines)
$first.Range.Text = $full

# Turn the sentinels into real paragraph breaks / tab characters.
$d.Content.Find.Execute([string]$PB, $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)
$d.Content.Find.Execute([string]$TB, $false, $false, $false, $false, $false, $true, 1, $false, "^t", 2)

# Groups (1-based paragraph indices) sharing the same paragraph formatting.
$groupA = 1..10     # spacing after=0, line=240/auto
$groupB = 11..16    # spacing after=0, line=240/auto, hanging indent 2127
$groupC = 17..18    # spacing after=0 only, hanging indent 2126

foreach ($i in $groupA) {
    $p = $d.Paragraphs($i)
    $p.SpaceAfter = 0
    $p.LineSpacingRule = 0
    $p.LineSpacing = 12
}

foreach ($i in $groupB) {
    $p = $d.Paragraphs($i)
    $p.SpaceAfter = 0
    $p.LineSpacingRule = 0
    $p.LineSpacing = 12
    $p.LeftIndent = 106.35
    $p.FirstLineIndent = -106.35
}

foreach ($i in $groupC) {
    $p = $d.Paragraphs($i)
    $p.SpaceAfter = 0
    $p.LeftIndent = 106.3
    $p.FirstLineIndent = -106.3
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
